# Apply the "Review first 10 problems" update to the LeetCode tags sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New "Palindrome" mini-table below the main tag table ----------------
$ws.Range("A24").Value = "Palindrome"
$ws.Range("A25").Value = 5

# --- Row 1 headers: new tag columns inserted after the original D1 -------
$ws.Range("E1").Value = "Binary Search"
$ws.Range("E2").Value = 5

$ws.Range("G1").Value = "strings"
$ws.Range("H1").Value = " Bits"
$ws.Range("G2").Value = "6,8，10"
$ws.Range("H2").Value = 7

$ws.Range("I1").Value = "DP"
$ws.Range("I2").Value = 10

$ws.Range("P1").Value = "Too Hard"
$ws.Range("P2").Value = 10

$ws.Range("J1").Value = "Greedy"
$ws.Range("J2").Value = 12

# --- Old far-right "ATTENTION！" column moves from R to O -----------------
$ws.Range("O1").Value = "ATTENTION！"
$ws.Range("O2").Value = 4
$ws.Range("R1").Value = $null
$ws.Range("R2").Value = $null

# --- View / selection tweaks ----------------------------------------------
[void]$ws.Range("J2").Select()

# --- Page setup (new pageSetup element appears in the diff) --------------
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
